$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns at G:I. This shifts the existing
# "mol/dL" formula column G->J, and the hemocue columns I,J,K,L -> L,M,N,O
# (their shared-string header references & values move with them automatically).
$ws.Range("G1:I1").EntireColumn.Insert()

# New header labels for the newly inserted / repurposed columns.
# Order matters: it controls the append order into the shared-string table.
$ws.Range("K1").Value = "c mol/dL"
$ws.Range("G1").Value = "cadj mol/dL"
$ws.Range("F1").Value = "cadj g/dL"
$ws.Range("H1").Value = "hpadj g/dL"
$ws.Range("I1").Value = "hpadj mol/dL"

# Raw measured values for the new "cadj g/dL" (F) and "hpadj g/dL" (H) columns
$ws.Range("F2").Value = -2.4374999999999947
$ws.Range("H2").Value = -2.145833333333328
$ws.Range("F3").Value = -4.75663333333333
$ws.Range("H3").Value = -1.889933333333329
$ws.Range("F4").Value = -9.278199999999995
$ws.Range("H4").Value = -2.6072999999999955
$ws.Range("F5").Value = -13.606366666666663
$ws.Range("H5").Value = -4.172933333333327
$ws.Range("F6").Value = -16.199466666666662
$ws.Range("H6").Value = 26.517333333333333
$ws.Range("F7").Value = 25.650833333333342
$ws.Range("H7").Value = 44.92183333333334
$ws.Range("F8").Value = 107.4237
$ws.Range("H8").Value = 32.22806666666667
$ws.Range("F9").Value = 153.0358
$ws.Range("H9").Value = 37.18603333333334
$ws.Range("F10").Value = 156.19546666666668
$ws.Range("H10").Value = 50.420733333333345
$ws.Range("F11").Value = 161.31936666666667
$ws.Range("H11").Value = 74.328
$ws.Range("F12").Value = 163.6575
$ws.Range("H12").Value = 115.03283333333336
$ws.Range("F13").Value = 162.93903333333333
$ws.Range("H13").Value = 129.84773333333334
$ws.Range("F14").Value = 172.5598
$ws.Range("H14").Value = 166.56436666666667
$ws.Range("F15").Value = 182.99896666666666
$ws.Range("H15").Value = 167.7452333333333
$ws.Range("F16").Value = 177.63153333333332
$ws.Range("H16").Value = 184.18200000000002
$ws.Range("F17").Value = 173.97833333333335
$ws.Range("H17").Value = 178.333
$ws.Range("F18").Value = 185.6371666666667
$ws.Range("H18").Value = 194.13766666666663
$ws.Range("F19").Value = 178.5043333333333
$ws.Range("H19").Value = 189.58816666666667
$ws.Range("F20").Value = 189.85066666666668
$ws.Range("H20").Value = 190.28866666666664
$ws.Range("F21").Value = 174.4886666666667
$ws.Range("H21").Value = 182.13499999999996
$ws.Range("F22").Value = 189.1683333333333
$ws.Range("H22").Value = 179.96049999999997
$ws.Range("F23").Value = 180.55633333333333
$ws.Range("H23").Value = 167.911
$ws.Range("F24").Value = 178.59016666666662
$ws.Range("H24").Value = 170.34066666666664
$ws.Range("F25").Value = 173.29066666666662
$ws.Range("H25").Value = 165.81199999999995

# Formulas, entered in the same order as the target shared-formula group ids:
# J (mol/dL, was G) and K (c mol/dL, new) filled down together from row 2 first ...
$ws.Range("J2:J25").Formula = "=D2/(1000*64000)"
$ws.Range("K2:K25").Formula = "=B2/(1000*64000)"

# ... then G/I are entered once at row 2 (standalone) and filled down from row 3
# (matches the shared-formula grouping in the target file: G2/I2 standalone, G3:G25/I3:I25 shared)
$ws.Range("G2").Formula = "=F2/(100*64000)"
$ws.Range("I2").Formula = "=H2/(100*64000)"
$ws.Range("G3:G25").Formula = "=F3/(100*64000)"
$ws.Range("I3:I25").Formula = "=H3/(100*64000)"

